$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header F1: "Wheat.Phenology.Stage", same style as the other header cells
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Wheat.Phenology.Stage"

# New observation rows (358-405): Wheat.Phenology.Stage values per simulation/date
$data = @(
    @{ Row=358; Name="Dookie2024CvGregory"; Date=45438; F=3 },
    @{ Row=359; Name="Dookie2024CvIllabo"; Date=45409; F=3 },
    @{ Row=360; Name="Dookie2024CvJanz"; Date=45439; F=3 },
    @{ Row=361; Name="Dookie2024CvKittyhawk"; Date=45409; F=3 },
    @{ Row=362; Name="Dookie2024CvMace"; Date=45439; F=3 },
    @{ Row=363; Name="Dookie2024CvMeering"; Date=45437; F=3 },
    @{ Row=364; Name="Dookie2024CvMowhawk"; Date=45410; F=3 },
    @{ Row=365; Name="Dookie2024CvOsprey"; Date=45416; F=3 },
    @{ Row=366; Name="Dookie2024CvRosella"; Date=45414; F=3 },
    @{ Row=367; Name="Dookie2024CvScepter"; Date=45439; F=3 },
    @{ Row=368; Name="Dookie2024CvSunmaster"; Date=45439; F=3 },
    @{ Row=369; Name="Dookie2024CvWedgetail"; Date=45410; F=3 },
    @{ Row=370; Name="Dookie2024CvWhistler"; Date=45413; F=3 },
    @{ Row=371; Name="Dookie2024CvWyalkatchem"; Date=45438; F=3 },
    @{ Row=372; Name="Dookie2024CvWylah"; Date=45414; F=3 },
    @{ Row=373; Name="Dookie2024CvYitpi"; Date=45438; F=3 },
    @{ Row=374; Name="Dookie2024CvGregory"; Date=45547; F=6 },
    @{ Row=375; Name="Dookie2024CvIllabo"; Date=45543; F=6 },
    @{ Row=376; Name="Dookie2024CvJanz"; Date=45545; F=6 },
    @{ Row=377; Name="Dookie2024CvKittyhawk"; Date=45547; F=6 },
    @{ Row=378; Name="Dookie2024CvMace"; Date=45541; F=6 },
    @{ Row=379; Name="Dookie2024CvMeering"; Date=45549; F=6 },
    @{ Row=380; Name="Dookie2024CvMowhawk"; Date=45538; F=6 },
    @{ Row=381; Name="Dookie2024CvOsprey"; Date=45545; F=6 },
    @{ Row=382; Name="Dookie2024CvRosella"; Date=45548; F=6 },
    @{ Row=383; Name="Dookie2024CvScepter"; Date=45542; F=6 },
    @{ Row=384; Name="Dookie2024CvSunmaster"; Date=45546; F=6 },
    @{ Row=385; Name="Dookie2024CvWedgetail"; Date=45544; F=6 },
    @{ Row=386; Name="Dookie2024CvWhistler"; Date=45539; F=6 },
    @{ Row=387; Name="Dookie2024CvWyalkatchem"; Date=45546; F=6 },
    @{ Row=388; Name="Dookie2024CvWylah"; Date=45546; F=6 },
    @{ Row=389; Name="Dookie2024CvYitpi"; Date=45546; F=6 },
    @{ Row=390; Name="Dookie2024CvGregory"; Date=45568; F=8 },
    @{ Row=391; Name="Dookie2024CvIllabo"; Date=45565; F=8 },
    @{ Row=392; Name="Dookie2024CvJanz"; Date=45566; F=8 },
    @{ Row=393; Name="Dookie2024CvKittyhawk"; Date=45568; F=8 },
    @{ Row=394; Name="Dookie2024CvMace"; Date=45560; F=8 },
    @{ Row=395; Name="Dookie2024CvMeering"; Date=45567; F=8 },
    @{ Row=396; Name="Dookie2024CvMowhawk"; Date=45560; F=8 },
    @{ Row=397; Name="Dookie2024CvOsprey"; Date=45566; F=8 },
    @{ Row=398; Name="Dookie2024CvRosella"; Date=45567; F=8 },
    @{ Row=399; Name="Dookie2024CvScepter"; Date=45561; F=8 },
    @{ Row=400; Name="Dookie2024CvSunmaster"; Date=45567; F=8 },
    @{ Row=401; Name="Dookie2024CvWedgetail"; Date=45566; F=8 },
    @{ Row=402; Name="Dookie2024CvWhistler"; Date=45560; F=8 },
    @{ Row=403; Name="Dookie2024CvWyalkatchem"; Date=45566; F=8 },
    @{ Row=404; Name="Dookie2024CvWylah"; Date=45567; F=8 },
    @{ Row=405; Name="Dookie2024CvYitpi"; Date=45567; F=8 }
)

# Apply the same formatting as existing data rows to the new rows in one shot
$ws.Range("A2").Copy()
$ws.Range("A358:A405").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B358:B405").PasteSpecial(-4122)

foreach ($d in $data) {
    $ws.Cells.Item($d.Row, 1).Value = $d.Name
    $ws.Cells.Item($d.Row, 2).Value = $d.Date
    $ws.Cells.Item($d.Row, 6).Value = $d.F
}
